# SyncLab_SoftEdgeEffect: give the "Source" rectangles an explicit accent1
# fill + a thin outline, bump their Soft Edge radius from a round 10pt to a
# decimal 30.23pt, and add a (previously implicit) Soft Edge effect to the
# "Destination" triangles (0pt on slide 1, 30.23pt on slide 2 to match the
# rectangles there).

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    $source = $s.Shapes.Item(1)
    $destination = $s.Shapes.Item(2)

    # Rectangle ("Source"): solid accent1 fill.
    $source.Fill.ForeColor.ObjectThemeColor = 5

    # Rectangle ("Source"): thin outline, #41719C, matching PowerPoint's
    # default "Format Shape > Line" values once a line color is applied.
    $source.Line.Visible = 1
    $source.Line.ForeColor.RGB = 10252609
    $source.Line.Weight = 1
    $source.Line.Style = 1
    $source.Line.DashStyle = 1
    $source.Line.CapStyle = 3
    $source.Line.JoinStyle = 3
    $source.Line.BeginArrowheadStyle = 1
    $source.Line.BeginArrowheadLength = 2
    $source.Line.BeginArrowheadWidth = 2
    $source.Line.EndArrowheadStyle = 1
    $source.Line.EndArrowheadLength = 2
    $source.Line.EndArrowheadWidth = 2

    # Rectangle ("Source"): Soft Edge radius 10pt -> 30.23pt (decimal).
    $source.SoftEdge.Radius = 30.23

    # Triangle ("Destination"): add an explicit Soft Edge effect.
    if ($i -eq 1) {
        $destination.SoftEdge.Radius = 0
    } else {
        $destination.SoftEdge.Radius = 30.23
    }
}
